$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.207.65"
$ws.Range("E2").Value = "  +2.31%  "

$ws.Range("D3").Value = "2.378.26"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("E5").Value = "  +7.28%  "

$ws.Range("D6").Formula = "'243.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.33%  "

$ws.Range("D7").Formula = "'77.19"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.35%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Formula = "'0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +27.62%  "

$ws.Range("E10").Value = "  +7.53%  "

$ws.Range("D11").Formula = "'57.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("D12").Formula = "'32.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +19.64%  "

$ws.Range("D13").Formula = "'7.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +21.38%  "

$ws.Range("E14").Value = "  +2.64%  "

$ws.Range("D15").Value = "2.733.22"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Formula = "'17.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.94%  "

$ws.Range("D17").Formula = "'0.931"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.81%  "

$ws.Range("D18").Value = "2.368.99"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "44.483.51"
$ws.Range("E19").Value = "  +2.96%  "

$ws.Range("E20").Value = "  +4.68%  "

$ws.Range("E21").Value = "  +6.74%  "

$ws.Range("D22").Formula = "'78.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.44%  "

$ws.Range("D23").Formula = "'257.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.62%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("E25").Value = "  +5.47%  "

$ws.Range("E26").Value = "  +0.61%  "

$ws.Range("D27").Formula = "'10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.61%  "

$ws.Range("D28").Formula = "'1.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +19.76%  "

$ws.Range("E29").Value = "  +2.12%  "

$ws.Range("D30").Formula = "'23.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.94%  "

$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("E32").Value = "  +1.84%  "

$ws.Range("E33").Value = "  +7.40%  "

$ws.Range("D34").Formula = "'5.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.76%  "

$ws.Range("D35").Formula = "'0.0765"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.05%  "

$ws.Range("E36").Value = "  +6.60%  "

$ws.Range("E37").Value = "  +5.59%  "

$ws.Range("D38").Formula = "'2.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.76%  "

$ws.Range("D40").Formula = "'0.0279"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.92%  "

$ws.Range("D41").Formula = "'9.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.46%  "

$ws.Range("D42").Formula = "'19.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.70%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Formula = "'0.197"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +18.63%  "

$ws.Range("E45").Value = "  +4.19%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Formula = "'1.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.24%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Formula = "'2.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +13.95%  "

$ws.Range("E48").Value = "  +6.11%  "

$ws.Range("D49").Formula = "'103.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("E50").Value = "  -1.17%  "

$ws.Range("D51").Value = "1.478.27"
$ws.Range("E51").Value = "  +2.39%  "

